# PROS-7621 - MARSRU - New OSA
# Updates KPI text (adds "Есть ЦА:" / "Нет ЦА:" prefixes for new OSA logic),
# updates a couple of existing KPI id rows, and appends four new KPI rows
# (54-57) for the "no central aisle" variants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update existing KPI display texts to reflect the "Есть ЦА" (has a
#    central aisle) scenario.
# ---------------------------------------------------------------------
$ws.Range("D35").Value = "Есть ЦА: Категория выстроена либо в единую линию, либо в две линии строго друг напротив друга (лицом друг к другу)"
$ws.Range("D37").Value = "Есть ЦА: Категория товаров для животных примыкает к ЦЕНТРАЛЬНОЙ АЛЛЕЕ и визуально доступна покупателям по ходу их движения без необходимости оборачиваться"
$ws.Range("D38").Value = "Есть ЦА: Категория товаров для животных примыкает к ПРОМО АЛЛЕЕ, находится дальше 5-ти метров от входа и визуально доступна покупателям по ходу их движения без необходимости оборачиваться"

# ---------------------------------------------------------------------
# 2. Renumber two existing KPI ids (rows 50 & 51).
# ---------------------------------------------------------------------
$ws.Range("B50").Value = 4697
$ws.Range("C50").Value = 4697

$ws.Range("B51").Value = 4698
$ws.Range("C51").Value = 4698

# Row heights for rows 50, 51 and 53 settle to 12.8 (from 13.8) after
# the edit, same as the rest of the "standard" (non-18.8) rows.
$ws.Rows.Item(50).RowHeight = 12.8
$ws.Rows.Item(51).RowHeight = 12.8
$ws.Rows.Item(53).RowHeight = 12.8

# ---------------------------------------------------------------------
# 3. Append four new KPI rows (54-57) for the "Нет ЦА" (no central
#    aisle) scenario. Copy formatting from row 53 (same column layout)
#    before filling in the new content.
# ---------------------------------------------------------------------
$ws.Range("A53:D53").Copy()
$ws.Range("A54:D57").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A54").Value = "MARS KPIs"
$ws.Range("B54").Value = 4601
$ws.Range("C54").Value = 4601
$ws.Range("D54").Value = "Нет ЦА: Категория тов. для животных примыкает или расположена в радиусе 5 м от центра выкладки приоритетной категории (1) (молочные прод, фрукты и овощи, хлебобулочные изд, кондитерские изд, мясн. изд. и рыба) таким образом, что видны блоки паучей Kitekat и Whiskas"

$ws.Range("A55").Value = "MARS KPIs"
$ws.Range("B55").Value = 4602
$ws.Range("C55").Value = 4602
$ws.Range("D55").Value = "Нет ЦА: Категория товаров для животных примыкает или расположена в радиусе 5 м от центра выкладки приоритетной категории (2) (консервы, соки, вода/газированные напитки, замороженные продукты), таким образом, что видны блоки паучей Kitekat  и Whiskas"

$ws.Range("A56").Value = "MARS KPIs"
$ws.Range("B56").Value = 4603
$ws.Range("C56").Value = 4603
$ws.Range("D56").Value = "Нет ЦА: Категория выстроена в единую линию единым блоком или образует внутренний угол"

$ws.Range("A57").Value = "MARS KPIs"
$ws.Range("B57").Value = 4604
$ws.Range("C57").Value = 4604
$ws.Range("D57").Value = "Нет ЦА: Категория располагается вне тупика и находится дальше 5 м от входа/выхода и кассовой зоны"

$ws.Rows.Item(54).RowHeight = 12.8
$ws.Rows.Item(55).RowHeight = 12.8
$ws.Rows.Item(56).RowHeight = 12.8
$ws.Rows.Item(57).RowHeight = 12.8

# ---------------------------------------------------------------------
# 4. Keep the cursor / selection on the last new row, matching the
#    author's final position when the edit was made.
# ---------------------------------------------------------------------
$ws.Range("A55").Select()
